$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; Col = "D"; Value = "328.61" }
    @{ Row = 2; Col = "E"; Value = "1.74%" }
    @{ Row = 3; Col = "D"; Value = "41.41" }
    @{ Row = 3; Col = "E"; Value = "5.00%" }
    @{ Row = 4; Col = "D"; Value = "5.618" }
    @{ Row = 4; Col = "E"; Value = "-4.06%" }
    @{ Row = 5; Col = "D"; Value = "0.08170" }
    @{ Row = 5; Col = "E"; Value = "2.03%" }
    @{ Row = 6; Col = "E"; Value = "1.70%" }
    @{ Row = 7; Col = "D"; Value = "8.737" }
    @{ Row = 7; Col = "E"; Value = "1.38%" }
    @{ Row = 8; Col = "D"; Value = "4.521" }
    @{ Row = 8; Col = "E"; Value = "-0.88%" }
    @{ Row = 9; Col = "D"; Value = "2.963" }
    @{ Row = 9; Col = "E"; Value = "0.35%" }
    @{ Row = 10; Col = "D"; Value = "0.9213" }
    @{ Row = 10; Col = "E"; Value = "-0.85%" }
    @{ Row = 11; Col = "D"; Value = "0.1276" }
    @{ Row = 11; Col = "E"; Value = "0.75%" }
    @{ Row = 12; Col = "D"; Value = "0.1959" }
    @{ Row = 12; Col = "E"; Value = "0.42%" }
    @{ Row = 13; Col = "D"; Value = "0.09401" }
    @{ Row = 13; Col = "E"; Value = "2.69%" }
    @{ Row = 14; Col = "D"; Value = "0.03813" }
    @{ Row = 14; Col = "E"; Value = "5.13%" }
    @{ Row = 15; Col = "D"; Value = "0.1060" }
    @{ Row = 15; Col = "E"; Value = "1.09%" }
    @{ Row = 16; Col = "D"; Value = "0.001304" }
    @{ Row = 16; Col = "E"; Value = "1.70%" }
    @{ Row = 17; Col = "D"; Value = "0.006283" }
    @{ Row = 17; Col = "E"; Value = "-1.73%" }
    @{ Row = 19; Col = "D"; Value = "3.442" }
    @{ Row = 19; Col = "E"; Value = "2.66%" }
    @{ Row = 20; Col = "E"; Value = "-1.12%" }
    @{ Row = 21; Col = "D"; Value = "8.312" }
    @{ Row = 21; Col = "E"; Value = "-4.54%" }
    @{ Row = 22; Col = "E"; Value = "1.17%" }
    @{ Row = 23; Col = "D"; Value = "0.2413" }
    @{ Row = 23; Col = "E"; Value = "-1.39%" }
    @{ Row = 24; Col = "D"; Value = "0.04418" }
    @{ Row = 24; Col = "E"; Value = "0.05%" }
    @{ Row = 25; Col = "D"; Value = "0.001261" }
    @{ Row = 25; Col = "E"; Value = "-0.34%" }
    @{ Row = 26; Col = "D"; Value = "0.004303" }
    @{ Row = 26; Col = "E"; Value = "-2.43%" }
    @{ Row = 27; Col = "E"; Value = "2.57%" }
    @{ Row = 39; Col = "D"; Value = "0.02779" }
    @{ Row = 39; Col = "E"; Value = "11.22%" }
    @{ Row = 40; Col = "D"; Value = "0.05421" }
    @{ Row = 40; Col = "E"; Value = "4.23%" }
    @{ Row = 41; Col = "D"; Value = "0.007680" }
    @{ Row = 41; Col = "E"; Value = "2.69%" }
    @{ Row = 42; Col = "E"; Value = "1.19%" }
    @{ Row = 43; Col = "D"; Value = "0.008980" }
    @{ Row = 43; Col = "E"; Value = "-6.58%" }
    @{ Row = 44; Col = "E"; Value = "0.64%" }
    @{ Row = 45; Col = "D"; Value = "0.01168" }
    @{ Row = 45; Col = "E"; Value = "5.31%" }
    @{ Row = 46; Col = "D"; Value = "0.00006647" }
    @{ Row = 46; Col = "E"; Value = "-1.53%" }
    @{ Row = 47; Col = "E"; Value = "0.13%" }
    @{ Row = 48; Col = "E"; Value = "6.64%" }
    @{ Row = 49; Col = "D"; Value = "0.002283" }
    @{ Row = 49; Col = "E"; Value = "-0.35%" }
    @{ Row = 50; Col = "D"; Value = "0.00002104" }
    @{ Row = 50; Col = "E"; Value = "0.13%" }
    @{ Row = 51; Col = "E"; Value = "0.13%" }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Col + $u.Row)
    $cell.Value = "'" + $u.Value
    $cell.Style = "Normal"
}
